$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 7 de Mayo de 2020 a las 17:04"

# --- Update numeric stats for countries whose ranking order does not change ---

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 1266442
$ws.Cells.Item(4, 3).Value = 3350
$ws.Cells.Item(4, 4).Value = 213138
$ws.Cells.Item(4, 5).Value = 978356
$ws.Cells.Item(4, 7).Value = 149
$ws.Cells.Item(4, 8).Value = 74948

# Row 12: Brasil
$ws.Cells.Item(12, 2).Value = 127389
$ws.Cells.Item(12, 3).Value = 778
$ws.Cells.Item(12, 5).Value = 67414
$ws.Cells.Item(12, 7).Value = 17
$ws.Cells.Item(12, 8).Value = 8605

# Row 15: Canada
$ws.Cells.Item(15, 2).Value = 63895
$ws.Cells.Item(15, 3).Value = 399
$ws.Cells.Item(15, 5).Value = 31444
$ws.Cells.Item(15, 7).Value = 48
$ws.Cells.Item(15, 8).Value = 4280

# Row 19: Paises Bajos
$ws.Cells.Item(19, 6).Value = 584

# Row 74: Azerbaiyan
$ws.Cells.Item(74, 2).Value = 2204
$ws.Cells.Item(74, 3).Value = 77
$ws.Cells.Item(74, 4).Value = 1551
$ws.Cells.Item(74, 5).Value = 625

# Row 79: Bulgaria
$ws.Cells.Item(79, 2).Value = 1829
$ws.Cells.Item(79, 3).Value = 51
$ws.Cells.Item(79, 5).Value = 1361

# --- Rows 142-144: updated data causes a re-sort (descending by "Casos totales"),
#     so Birmania (updated) now ranks above Santo Tome y Principe and Republica del Chad ---

# Row 142 becomes Birmania
$ws.Cells.Item(142, 1).Value = "Birmania"
$ws.Cells.Item(142, 2).Value = 176
$ws.Cells.Item(142, 3).Value = 15
$ws.Cells.Item(142, 4).Value = 62
$ws.Cells.Item(142, 5).Value = 108
$ws.Cells.Item(142, 6).Value = 0
$ws.Cells.Item(142, 7).Value = 0
$ws.Cells.Item(142, 8).Value = 6

# Row 143 becomes Santo Tome y Principe
$ws.Cells.Item(143, 1).Value = "Santo Tome y Principe"
$ws.Cells.Item(143, 2).Value = 174
$ws.Cells.Item(143, 3).Value = 0
$ws.Cells.Item(143, 4).Value = 4
$ws.Cells.Item(143, 5).Value = 167
$ws.Cells.Item(143, 6).Value = 0
$ws.Cells.Item(143, 7).Value = 0
$ws.Cells.Item(143, 8).Value = 3

# Row 144 becomes Republica del Chad
$ws.Cells.Item(144, 1).Value = "Republica del Chad"
$ws.Cells.Item(144, 2).Value = 170
$ws.Cells.Item(144, 3).Value = 0
$ws.Cells.Item(144, 4).Value = 43
$ws.Cells.Item(144, 5).Value = 110
$ws.Cells.Item(144, 6).Value = 0
$ws.Cells.Item(144, 7).Value = 0
$ws.Cells.Item(144, 8).Value = 17

# --- Rows 205-206: Montserrat and Seychelles swap positions (tied totals, data updated) ---

# Row 205 becomes Seychelles
$ws.Cells.Item(205, 1).Value = "Seychelles"
$ws.Cells.Item(205, 2).Value = 11
$ws.Cells.Item(205, 3).Value = 0
$ws.Cells.Item(205, 4).Value = 8
$ws.Cells.Item(205, 5).Value = 3
$ws.Cells.Item(205, 6).Value = 0
$ws.Cells.Item(205, 7).Value = 0
$ws.Cells.Item(205, 8).Value = 0

# Row 206 becomes Montserrat
$ws.Cells.Item(206, 1).Value = "Montserrat"
$ws.Cells.Item(206, 2).Value = 11
$ws.Cells.Item(206, 3).Value = 0
$ws.Cells.Item(206, 4).Value = 7
$ws.Cells.Item(206, 5).Value = 3
$ws.Cells.Item(206, 6).Value = 1
$ws.Cells.Item(206, 7).Value = 0
$ws.Cells.Item(206, 8).Value = 1
